{"js": "// Word template now needs distinct even/default/first page headers & footers\n// (previously only a \"first page\" header/footer existed because of titlePg),\n// and the \"first page\" footer gains a new \"Victim's Attorney\" line item.\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst sec = sections.items[0];\n\n// Touching the Primary (default) header/footer on a section that only has a\n// \"first page\" header/footer forces Word to materialize the full even /\n// default / first trio of header & footer parts (matching the headerReference\n// / footerReference w:type=\"even|default|first\" set added by the diff).\nconst primaryHeader = sec.getHeader(\"Primary\");\nconst primaryFooter = sec.getFooter(\"Primary\");\nprimaryHeader.clear();\nprimaryFooter.clear();\nawait context.sync();\n\n// The existing content (title-page header/footer) now lives under the\n// \"FirstPage\" header/footer. Append the new \"Victim's Attorney\" text to the\n// paragraph that already holds the \"Prosecutor's Office / ... / County Jail\"\n// service-of-copies line.\nconst firstFooter = sec.getFooter(\"FirstPage\");\nconst paragraphs = firstFooter.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Prosecutor\\u2019s Office\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  const appended = targetParagraph.insertText(\n    \" Victim\\u2019s Attorney (if applicable): PS   OS   EM\",\n    Word.InsertLocation.end\n  );\n  appended.font.name = \"Palatino Linotype\";\n  appended.font.size = 8;\n  await context.sync();\n}\n", "ps1": "# Word template now needs distinct even/default/first page headers & footers\n# (previously only a \"first page\" header/footer existed because of titlePg),\n# and the \"first page\" footer gains a new \"Victim's Attorney\" line item.\n\n$d = $word.ActiveDocument\n$sec = $d.Sections(1)\n\n# wdHeaderFooterIndex: 1 = Primary (default), 2 = FirstPage, 3 = EvenPages\n$primaryHeader = $sec.Headers(1)\n$primaryFooter = $sec.Footers(1)\n\n# Touching the Primary (default) header/footer on a section that only has a\n# \"first page\" header/footer forces Word to materialize the full even /\n# default / first trio of header & footer parts (matching the\n# headerReference / footerReference w:type=\"even|default|first\" set added by\n# the diff).\n$primaryHeader.Range.Text = \"\"\n$primaryFooter.Range.Text = \"\"\n\n# The existing content (title-page header/footer) now lives under the\n# \"FirstPage\" header/footer (index 2). Find the service-of-copies line\n# (\"Prosecutor's Office ... County Jail ...\") and append the new\n# \"Victim's Attorney\" sentence right after it, keeping the same formatting.\n$firstFooter = $sec.Footers(2)\n\n$find = $firstFooter.Range.Find\n$find.ClearFormatting()\n$found = $find.Execute(\"County Jail: PS   EM;\")\n\nif ($found) {\n    $insertPoint = $firstFooter.Range.Duplicate()\n    $insertPoint.SetRange($find.Parent.End, $find.Parent.End)\n    $insertPoint.Font.Name = \"Palatino Linotype\"\n    $insertPoint.Font.Size = 8\n    $insertPoint.InsertAfter(\" Victim\u2019s Attorney (if applicable): PS   OS   EM\")\n}\n"}
